$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "#! END_ROW true"
$ws.Range("E5").Value = "#! END_ROW true"
$ws.Range("E8").Value = "#! END_ROW true"

$ws.Range("E8").Select() | Out-Null
